$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as plain text in this workbook (e.g.
# "79.278.74" / "3.179.33" are not valid numbers), so force a Text number format
# on every Price cell we are about to rewrite before assigning the new value.
# This keeps Excel from "helpfully" reinterpreting them as numeric values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "79.278.74"
$ws.Range("E2").Value = "  +3.84%  "

$ws.Range("D3").Value = "3.179.33"
$ws.Range("E3").Value = "  +4.71%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "205.63"
$ws.Range("E5").Value = "  +3.78%  "

$ws.Range("D6").Value = "633.61"
$ws.Range("E6").Value = "  +2.56%  "

$ws.Range("E8").Value = "  +17.40%  "

$ws.Range("D9").Value = "0.598"
$ws.Range("E9").Value = "  +9.53%  "

$ws.Range("D10").Value = "3.174.52"
$ws.Range("E10").Value = "  +4.67%  "

$ws.Range("D11").Value = "0.598"
$ws.Range("E11").Value = "  +37.51%  "

$ws.Range("D12").Value = "0.0000252"
$ws.Range("E12").Value = "  +31.54%  "

$ws.Range("D13").Value = "0.165"
$ws.Range("E13").Value = "  +3.21%  "

$ws.Range("D14").Value = "5.37"
$ws.Range("E14").Value = "  +2.09%  "

$ws.Range("D15").Value = "3.759.38"

$ws.Range("D16").Value = "31.89"
$ws.Range("E16").Value = "  +10.87%  "

$ws.Range("D17").Value = "79.211.69"
$ws.Range("E17").Value = "  +3.81%  "

$ws.Range("D18").Value = "3.168.54"
$ws.Range("E18").Value = "  +4.49%  "

$ws.Range("D19").Value = "14.50"
$ws.Range("E19").Value = "  +7.71%  "

$ws.Range("D20").Value = "9.37"
$ws.Range("E20").Value = "  +4.80%  "

$ws.Range("D21").Value = "439.79"
$ws.Range("E21").Value = "  +15.94%  "

$ws.Range("E22").Value = "  +25.58%  "

$ws.Range("E23").Value = "  +20.43%  "

$ws.Range("D24").Value = "4.83"
$ws.Range("E24").Value = "  +11.47%  "

$ws.Range("D25").Value = "77.17"
$ws.Range("E25").Value = "  +6.09%  "

$ws.Range("D26").Value = "10.84"
$ws.Range("E26").Value = "  +12.00%  "

$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("E28").Value = "  +13.41%  "

$ws.Range("D29").Value = "9.08"

$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("E31").Value = "  +9.44%  "

$ws.Range("D32").Value = "538.12"
$ws.Range("E32").Value = "  +9.86%  "

$ws.Range("D33").Value = "0.155"
$ws.Range("E33").Value = "  +32.21%  "

$ws.Range("D34").Value = "2.03"
$ws.Range("E34").Value = "  +5.78%  "

$ws.Range("D35").Value = "22.93"
$ws.Range("E35").Value = "  +11.91%  "

$ws.Range("E36").Value = "  +16.11%  "

$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("E38").Value = "  +7.12%  "

$ws.Range("D39").Value = "163.56"
$ws.Range("E39").Value = "  +1.04%  "

$ws.Range("D40").Value = "20.01"
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "191.95"
$ws.Range("E41").Value = "  +0.91%  "

$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").Value = "5.52"
$ws.Range("E43").Value = "  +9.17%  "

$ws.Range("D44").Value = "1.80"
$ws.Range("E44").Value = "  +10.35%  "

$ws.Range("D45").Value = "0.796"
$ws.Range("E45").Value = "  +1.04%  "

$ws.Range("D46").Value = "2.65"
$ws.Range("E46").Value = "  +10.65%  "

$ws.Range("E47").Value = "  +5.63%  "

$ws.Range("D48").Value = "43.12"
$ws.Range("E48").Value = "  +2.79%  "

$ws.Range("D49").Value = "25.72"
$ws.Range("E49").Value = "  +16.66%  "

$ws.Range("D50").Value = "0.638"
$ws.Range("E50").Value = "  +6.63%  "

$ws.Range("D51").Value = "4.20"
$ws.Range("E51").Value = "  +8.65%  "
